# "Generate Report for Archive"
# The localization status for the two handed-back files moved on from
# "Ready for handoff" to "In Translation" on every sheet that surfaces a
# Status column (Overview's per-locale summary columns, and each locale's
# own Status column), and the Status-ish columns are re-auto-fit now that
# the label is shorter.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) show the status of
#     each handed-back file in rows 2-3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E1:F3").EntireColumn.AutoFit()

# --- zh-cn sheet: column C ("Status") rows 2-3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C1:C3").EntireColumn.AutoFit()

# --- de-de sheet: column C ("Status") rows 2-3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C1:C3").EntireColumn.AutoFit()
